$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 7217.1
$ws.Range("I106").Value = 7546.4443
$ws.Range("K106").Value = 7546.4443
$ws.Range("M106").Value = -6915.4443

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1738.2858
$ws.Range("I125").Value = 1644
$ws.Range("J125").Value = 1809
$ws.Range("K125").Value = 14796
$ws.Range("L125").Value = 16281
$ws.Range("M125").Value = -12336
$ws.Range("N125").Value = -21201

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 10107461
$ws.Range("I132").Value = 11115811
$ws.Range("K132").Value = 33347433
$ws.Range("M132").Value = -33344903

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 47619908
$ws.Range("I135").Value = 379.73334
$ws.Range("K135").Value = 3417.60006
$ws.Range("M135").Value = -882.6000599999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1672.5857
$ws.Range("I138").Value = 1117.8572
$ws.Range("J138").Value = 1910.3265
$ws.Range("K138").Value = 3353.5716
$ws.Range("L138").Value = 5730.979499999999
$ws.Range("M138").Value = 1786.4284
$ws.Range("N138").Value = -16010.9795

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 607.7778
$ws.Range("I141").Value = 558.75
$ws.Range("K141").Value = 1676.25
$ws.Range("M141").Value = 3503.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1174.8
$ws.Range("I61").Value = 1174.8
$ws.Range("K61").Value = 1174.8
$ws.Range("M61").Value = -962.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1052.2
$ws.Range("I74").Value = 705.43475
$ws.Range("J74").Value = 1716.8334
$ws.Range("K74").Value = 705.43475
$ws.Range("L74").Value = 1716.8334
$ws.Range("M74").Value = 168.56525
$ws.Range("N74").Value = -3464.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1052.2
$ws.Range("I77").Value = 705.43475
$ws.Range("J77").Value = 1716.8334
$ws.Range("K77").Value = 3527.17375
$ws.Range("L77").Value = 8584.166999999999
$ws.Range("M77").Value = 840.8262500000001
$ws.Range("N77").Value = -17320.167

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2289
$ws.Range("I132").Value = 1880.2667
$ws.Range("K132").Value = 5640.800099999999
$ws.Range("M132").Value = -3110.800099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1174.8
$ws.Range("I136").Value = 1174.8
$ws.Range("K136").Value = 3524.4
$ws.Range("M136").Value = -974.3999999999996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8761.5
$ws.Range("I134").Value = 1340.5
$ws.Range("K134").Value = 4021.5
$ws.Range("M134").Value = -1486.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1230.0588
$ws.Range("I31").Value = 1204.66
$ws.Range("K31").Value = 1204.66
$ws.Range("M31").Value = -909.6600000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1230.0588
$ws.Range("I34").Value = 1204.66
$ws.Range("K34").Value = 1204.66
$ws.Range("M34").Value = -1002.66

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1598.3334
$ws.Range("I58").Value = 1435
$ws.Range("K58").Value = 1435
$ws.Range("M58").Value = -1232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2061.3333
$ws.Range("I132").Value = 1000.75
$ws.Range("K132").Value = 3002.25
$ws.Range("M132").Value = -472.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1327.9395
$ws.Range("I134").Value = 1223.28
$ws.Range("J134").Value = 1655
$ws.Range("K134").Value = 3669.84
$ws.Range("L134").Value = 4965
$ws.Range("M134").Value = -1134.84
$ws.Range("N134").Value = -10035

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1598.3334
$ws.Range("I136").Value = 1435
$ws.Range("K136").Value = 4305
$ws.Range("M136").Value = -1755

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4540.8887
$ws.Range("J64").Value = 4691.5293
$ws.Range("L64").Value = 14074.5879
$ws.Range("N64").Value = -14614.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 4540.8887
$ws.Range("J67").Value = 4691.5293
$ws.Range("L67").Value = 14074.5879
$ws.Range("N67").Value = -15946.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1871.4286
$ws.Range("I81").Value = 1025
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 3075
$ws.Range("L81").Value = 9000
$ws.Range("M81").Value = -1952
$ws.Range("N81").Value = -11246

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 1871.4286
$ws.Range("I84").Value = 1025
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 9225
$ws.Range("L84").Value = 27000
$ws.Range("M84").Value = -3609
$ws.Range("N84").Value = -38232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 300
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 2370.4285
$ws.Range("J87").Value = 3766.3333
$ws.Range("L87").Value = 11298.9999
$ws.Range("N87").Value = -13794.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 9000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 9000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 27000
$ws.Range("N88").Value = -27856
$ws.Range("M88").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 300
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 2370.4285
$ws.Range("J90").Value = 3766.3333
$ws.Range("L90").Value = 33896.9997
$ws.Range("N90").Value = -46376.9997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 9000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 9000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 27000
$ws.Range("N91").Value = -29964
$ws.Range("M91").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3884.8
$ws.Range("I94").Value = 3174.6667
$ws.Range("J94").Value = 4950
$ws.Range("K94").Value = 9524.000100000001
$ws.Range("L94").Value = 14850
$ws.Range("M94").Value = -8848.000100000001
$ws.Range("N94").Value = -16202

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 500
$ws.Range("J98").Value = 500
$ws.Range("L98").Value = 1500
$ws.Range("N98").Value = -4496

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 2114
$ws.Range("I99").Value = 320
$ws.Range("J99").Value = 2512.6667
$ws.Range("K99").Value = 960
$ws.Range("L99").Value = 7538.000100000001
$ws.Range("M99").Value = 1286
$ws.Range("N99").Value = -12030.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 3267.2
$ws.Range("J100").Value = 3267.2
$ws.Range("L100").Value = 9801.599999999999
$ws.Range("N100").Value = -11423.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 2392.3635
$ws.Range("J106").Value = 2392.3635
$ws.Range("L106").Value = 7177.0905
$ws.Range("N106").Value = -9069.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 17859914
$ws.Range("I131").Value = 100000390
$ws.Range("J131").Value = 3288.652
$ws.Range("K131").Value = 300001170
$ws.Range("L131").Value = 9865.956
$ws.Range("M131").Value = -299996130
$ws.Range("N131").Value = -19945.956

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1734.8064
$ws.Range("I132").Value = 1412.6522
$ws.Range("K132").Value = 4237.9566
$ws.Range("M132").Value = -1707.9566

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 29447.334
$ws.Range("I132").Value = 1061.762
$ws.Range("J132").Value = 69187.13
$ws.Range("K132").Value = 3185.286
$ws.Range("L132").Value = 207561.39
$ws.Range("M132").Value = -655.2860000000001
$ws.Range("N132").Value = -212621.39

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1052.5333
$ws.Range("I136").Value = 1004.4815
$ws.Range("J136").Value = 1485
$ws.Range("K136").Value = 3013.4445
$ws.Range("L136").Value = 4455
$ws.Range("M136").Value = -463.4445000000001
$ws.Range("N136").Value = -9555

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2191.2
$ws.Range("I132").Value = 1879.2222
$ws.Range("K132").Value = 5637.6666
$ws.Range("M132").Value = -3107.6666
